# Apply "5 and 10 mile radius" updates to the national transitions-rule
# summary tables workbook. Adds two new columns (Within 5 miles / Within
# 10 miles of HFC production facility) to both the "Means" and
# "Standard Deviations" sheets, and updates the existing Total Cancer
# Risk / Total Respiratory rows whose 1/3/overall/rural values shifted
# as part of the same re-run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Means"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cells for the two added distance bands
$ws1.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$ws1.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

# Row 2: % White
$ws1.Cells.Item(2, 6).Value = 59
$ws1.Cells.Item(2, 7).Value = 52

# Row 3: % Black or African American
$ws1.Cells.Item(3, 6).Value = 13
$ws1.Cells.Item(3, 7).Value = 22

# Row 4: % Other
$ws1.Cells.Item(4, 6).Value = 27
$ws1.Cells.Item(4, 7).Value = 26

# Row 5: % Hispanic
$ws1.Cells.Item(5, 6).Value = 41
$ws1.Cells.Item(5, 7).Value = 33

# Row 6: Median Income [1,000 2019$]
$ws1.Cells.Item(6, 6).Value = 67
$ws1.Cells.Item(6, 7).Value = 69

# Row 7: % Below Poverty Line
$ws1.Cells.Item(7, 6).Value = 8.3
$ws1.Cells.Item(7, 7).Value = 7.7

# Row 8: % Below Half the Poverty Line
$ws1.Cells.Item(8, 6).Value = 6.6
$ws1.Cells.Item(8, 7).Value = 6.3

# Row 9: Total Cancer Risk (per million) -- existing columns updated too
$ws1.Cells.Item(9, 2).Value = 29
$ws1.Cells.Item(9, 3).Value = 26
$ws1.Cells.Item(9, 4).Value = 44
$ws1.Cells.Item(9, 5).Value = 38
$ws1.Cells.Item(9, 6).Value = 37
$ws1.Cells.Item(9, 7).Value = 36

# Row 10: Total Respiratory (hazard quotient) -- existing columns updated too
$ws1.Cells.Item(10, 2).Value = 0.37
$ws1.Cells.Item(10, 3).Value = 0.32
$ws1.Cells.Item(10, 4).Value = 0.44
$ws1.Cells.Item(10, 5).Value = 0.44
$ws1.Cells.Item(10, 6).Value = 0.43
$ws1.Cells.Item(10, 7).Value = 0.42

# ---------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New header cells for the two added distance bands
$ws2.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$ws2.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

# Row 2: % White
$ws2.Cells.Item(2, 6).Value = 24
$ws2.Cells.Item(2, 7).Value = 31

# Row 3: % Black or African American
$ws2.Cells.Item(3, 6).Value = 23
$ws2.Cells.Item(3, 7).Value = 36

# Row 4: % Other
$ws2.Cells.Item(4, 6).Value = 20
$ws2.Cells.Item(4, 7).Value = 23

# Row 5: % Hispanic
$ws2.Cells.Item(5, 6).Value = 26
$ws2.Cells.Item(5, 7).Value = 29

# Row 6: Median Income [1,000 2019$]
$ws2.Cells.Item(6, 6).Value = 29
$ws2.Cells.Item(6, 7).Value = 29

# Row 7: % Below Poverty Line
$ws2.Cells.Item(7, 6).Value = 9.1
$ws2.Cells.Item(7, 7).Value = 9.2

# Row 8: % Below Half the Poverty Line
$ws2.Cells.Item(8, 6).Value = 8.4
$ws2.Cells.Item(8, 7).Value = 8.6

# Row 9: Total Cancer Risk (per million) -- existing columns updated too
$ws2.Cells.Item(9, 2).Value = 10
$ws2.Cells.Item(9, 3).Value = 8.6
$ws2.Cells.Item(9, 4).Value = 25
$ws2.Cells.Item(9, 5).Value = 20
$ws2.Cells.Item(9, 6).Value = 17
$ws2.Cells.Item(9, 7).Value = 13

# Row 10: Total Respiratory (hazard quotient) -- existing columns updated too
$ws2.Cells.Item(10, 2).Value = 0.14
$ws2.Cells.Item(10, 3).Value = 0.14
$ws2.Cells.Item(10, 4).Value = 0.094
$ws2.Cells.Item(10, 5).Value = 0.082
$ws2.Cells.Item(10, 6).Value = 0.084
$ws2.Cells.Item(10, 7).Value = 0.075

Write-Output "Applied 5mi/10mi radius updates to Means and Standard Deviations sheets."
